$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting existing rows 2,3,5 down to 3,4,6
$ws.Rows("2:2").Insert()

# Update comment text in row 1
$ws.Range("A1").Value = "Put the path to images under the corresponding title."

# Set new row 2 comment text and merge it like the other comment rows
$ws.Range("A2").Value = "The files in each row must correspond to different polarizations of same sample."
$ws.Range("A2:F2").Merge()
